$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "28.190.39"
$ws.Cells.Item(2,5).Value = "  -2.85%  "

# Row 3
$ws.Cells.Item(3,4).Value = "1.930.21"
$ws.Cells.Item(3,5).Value = "  -1.50%  "

# Row 4
$ws.Cells.Item(4,4).Value = "'1.011"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "  +0.63%  "

# Row 5
$ws.Cells.Item(5,4).Value = "'321.72"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -1.48%  "

# Row 6
$ws.Cells.Item(6,4).Value = "'1.011"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  +0.72%  "

# Row 7
$ws.Cells.Item(7,4).Value = "'0.4732"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "  -4.67%  "

# Row 8
$ws.Cells.Item(8,4).Value = "'0.4051"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "  -3.40%  "

# Row 9
$ws.Cells.Item(9,4).Value = "'53.40"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "  +1.36%  "

# Row 10
$ws.Cells.Item(10,4).Value = "'0.08514"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "  -7.37%  "

# Row 11
$ws.Cells.Item(11,4).Value = "'1.050"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  -4.12%  "

# Row 12
$ws.Cells.Item(12,4).Value = "'22.22"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "  -2.40%  "

# Row 13
$ws.Cells.Item(13,2).Value = "WrappedEther"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13,4).Value = "1.978.35"
$ws.Cells.Item(13,5).Value = "  +0.45%  "

# Row 14
$ws.Cells.Item(14,2).Value = "Chainlink"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14,4).Value = "'7.518"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  -3.87%  "

# Row 15
$ws.Cells.Item(15,4).Value = "'6.114"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "  -4.98%  "

# Row 16
$ws.Cells.Item(16,4).Value = "'1.015"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "  +0.84%  "

# Row 17
$ws.Cells.Item(17,4).Value = "'89.86"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "  -1.53%  "

# Row 18
$ws.Cells.Item(18,5).Value = "  -2.30%  "

# Row 19
$ws.Cells.Item(19,4).Value = "'0.06608"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -1.16%  "

# Row 20
$ws.Cells.Item(20,5).Value = "  -5.39%  "

# Row 21
$ws.Cells.Item(21,4).Value = "'1.011"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +0.76%  "

# Row 22
$ws.Cells.Item(22,4).Value = "'5.778"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -2.60%  "

# Row 23
$ws.Cells.Item(23,4).Value = "28.206.79"
$ws.Cells.Item(23,5).Value = "  -2.92%  "

# Row 24
$ws.Cells.Item(24,5).Value = "  -4.87%  "

# Row 25
$ws.Cells.Item(25,4).Value = "'2.306"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  +1.65%  "

# Row 26
$ws.Cells.Item(26,4).Value = "2.128.49"
$ws.Cells.Item(26,5).Value = "  -3.64%  "

# Row 27
$ws.Cells.Item(27,4).Value = "'155.01"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  -0.35%  "

# Row 28
$ws.Cells.Item(28,4).Value = "'20.17"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "  -1.88%  "

# Row 29
$ws.Cells.Item(29,4).Value = "'2.162"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -3.74%  "

# Row 30
$ws.Cells.Item(30,4).Value = "'5.753"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "  -8.06%  "

# Row 31
$ws.Cells.Item(31,4).Value = "'123.81"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = "  -1.78%  "

# Row 32
$ws.Cells.Item(32,4).Value = "'0.9795"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = "  -5.74%  "

# Row 33
$ws.Cells.Item(33,4).Value = "'0.09597"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  -2.23%  "

# Row 34
$ws.Cells.Item(34,4).Value = "'1.446"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -4.92%  "

# Row 35
$ws.Cells.Item(35,5).Value = "  -0.22%  "

# Row 36
$ws.Cells.Item(36,4).Value = "'5.586"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  -3.95%  "

# Row 37
$ws.Cells.Item(37,4).Value = "'9.267"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "  +3.38%  "

# Row 38
$ws.Cells.Item(38,4).Value = "'0.02321"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  -3.99%  "

# Row 39
$ws.Cells.Item(39,4).Value = "'0.06172"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  -2.58%  "

# Row 40
$ws.Cells.Item(40,4).Value = "'1.238"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -6.31%  "

# Row 41
$ws.Cells.Item(41,4).Value = "'0.6186"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "  -3.75%  "

# Row 42
$ws.Cells.Item(42,4).Value = "'11.12"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "  -2.18%  "

# Row 43
$ws.Cells.Item(43,5).Value = "  +0.77%  "

# Row 44
$ws.Cells.Item(44,4).Value = "'0.1908"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "  -3.14%  "

# Row 45
$ws.Cells.Item(45,4).Value = "'1.321"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  -2.66%  "

# Row 46
$ws.Cells.Item(46,4).Value = "'0.5901"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  -4.93%  "

# Row 47
$ws.Cells.Item(47,4).Value = "'12.85"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -3.55%  "

# Row 48
$ws.Cells.Item(48,4).Value = "'2.043"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "  -6.71%  "

# Row 49
$ws.Cells.Item(49,4).Value = "'3.395"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "  -1.95%  "

# Row 50
$ws.Cells.Item(50,4).Value = "'0.06767"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  -3.32%  "

# Row 51
$ws.Cells.Item(51,2).Value = "EOS"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(51,4).Value = "'1.087"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -2.24%  "
